$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7495
$ws.Range("J51").Value = 7495
$ws.Range("L51").Value = 7495
$ws.Range("N51").Value = -8463

$ws.Range("H92").Value = 1366.7742
$ws.Range("I92").Value = 1407.2333
$ws.Range("K92").Value = 1407.2333
$ws.Range("M92").Value = -159.2333000000001

$ws.Range("H115").Value = 969.5833
$ws.Range("I115").Value = 563.4
$ws.Range("K115").Value = 1690.2
$ws.Range("M115").Value = -123.1999999999998

$ws.Range("H137").Value = 3292.475
$ws.Range("J137").Value = 3920
$ws.Range("L137").Value = 11760
$ws.Range("N137").Value = -16860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3588.602
$ws.Range("I32").Value = 1986.6267
$ws.Range("K32").Value = 1986.6267
$ws.Range("M32").Value = -1699.6267

$ws.Range("H61").Value = 3654.8
$ws.Range("I61").Value = 3717.7874
$ws.Range("J61").Value = 3284.75
$ws.Range("K61").Value = 3717.7874
$ws.Range("L61").Value = 3284.75
$ws.Range("M61").Value = -3505.7874
$ws.Range("N61").Value = -3708.75

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 1304133
$ws.Range("I122").Value = 2193.5454
$ws.Range("K122").Value = 6580.6362
$ws.Range("M122").Value = -4130.6362

$ws.Range("H132").Value = 5314.879
$ws.Range("I132").Value = 5340.9287
$ws.Range("J132").Value = 5169
$ws.Range("K132").Value = 16022.7861
$ws.Range("L132").Value = 15507
$ws.Range("M132").Value = -13492.7861
$ws.Range("N132").Value = -20567

$ws.Range("H136").Value = 3654.8
$ws.Range("I136").Value = 3717.7874
$ws.Range("J136").Value = 3284.75
$ws.Range("K136").Value = 11153.3622
$ws.Range("L136").Value = 9854.25
$ws.Range("M136").Value = -8603.3622
$ws.Range("N136").Value = -14954.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2018040.4
$ws.Range("I105").Value = 2085268.1
$ws.Range("K105").Value = 2085268.1
$ws.Range("M105").Value = -2083521.1

$ws.Range("H134").Value = 10009.793
$ws.Range("I134").Value = 9262.041999999999
$ws.Range("K134").Value = 27786.126
$ws.Range("M134").Value = -25251.126

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 659.4
$ws.Range("I22").Value = 449.25
$ws.Range("K22").Value = 449.25
$ws.Range("M22").Value = -99.25

$ws.Range("H132").Value = 1905.8723
$ws.Range("I132").Value = 1782.4524
$ws.Range("K132").Value = 5347.357199999999
$ws.Range("M132").Value = -2817.357199999999

$ws.Range("H134").Value = 24504.568
$ws.Range("I134").Value = 29926.527
$ws.Range("J134").Value = 11491.866
$ws.Range("K134").Value = 89779.58099999999
$ws.Range("L134").Value = 34475.598
$ws.Range("M134").Value = -87244.58099999999
$ws.Range("N134").Value = -39545.598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 125750
$ws.Range("I9").Value = 1000000
$ws.Range("J9").Value = 857.1429000000001
$ws.Range("K9").Value = 3000000
$ws.Range("L9").Value = 2571.4287
$ws.Range("M9").Value = -2999776
$ws.Range("N9").Value = -3019.4287

$ws.Range("H56").Value = 16672961
$ws.Range("I56").Value = 16672961
$ws.Range("K56").Value = 16672961
$ws.Range("M56").Value = -16672431

$ws.Range("H64").Value = 356
$ws.Range("I64").Value = 356
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1068
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -798
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 356
$ws.Range("I67").Value = 356
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1068
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -132
$ws.Range("N67").ClearContents()

$ws.Range("H104").Value = 2200
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2200
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 6600
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -11842

$ws.Range("H131").Value = 10686178
$ws.Range("I131").Value = 5210032
$ws.Range("K131").Value = 15630096
$ws.Range("M131").Value = -15625056

$ws.Range("H137").Value = 2870.9333
$ws.Range("I137").Value = 1706.7
$ws.Range("J137").Value = 5199.4
$ws.Range("K137").Value = 5120.1
$ws.Range("L137").Value = 15598.2
$ws.Range("M137").Value = -20.10000000000036
$ws.Range("N137").Value = -25798.2

$ws.Range("H139").Value = 1781.7142
$ws.Range("I139").Value = 1098
$ws.Range("J139").Value = 2693.3333
$ws.Range("K139").Value = 3294
$ws.Range("L139").Value = 8079.999899999999
$ws.Range("M139").Value = 1846
$ws.Range("N139").Value = -18359.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7813.892
$ws.Range("I132").Value = 5550.625
$ws.Range("K132").Value = 16651.875
$ws.Range("M132").Value = -14121.875

$ws.Range("H134").Value = 35310.4
$ws.Range("J134").Value = 35310.4
$ws.Range("L134").Value = 105931.2
$ws.Range("N134").Value = -111001.2

$ws.Range("H136").Value = 38173.4
$ws.Range("J136").Value = 38173.4
$ws.Range("L136").Value = 114520.2
$ws.Range("N136").Value = -119620.2

$ws.Range("H138").Value = 68604.25
$ws.Range("J138").Value = 67208.5
$ws.Range("L138").Value = 67208.5
$ws.Range("N138").Value = -77488.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1143.7059
$ws.Range("J16").Value = 1004
$ws.Range("L16").Value = 1004
$ws.Range("N16").Value = -1344

$ws.Range("H61").Value = 7410005
$ws.Range("I61").Value = 10103319
$ws.Range("J61").Value = 3392.25
$ws.Range("K61").Value = 10103319
$ws.Range("L61").Value = 3392.25
$ws.Range("M61").Value = -10103117
$ws.Range("N61").Value = -3796.25

$ws.Range("H68").Value = 2842.7144
$ws.Range("I68").Value = 2983.1667
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2983.1667
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -2234.1667
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 2842.7144
$ws.Range("I71").Value = 2983.1667
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 14915.8335
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -11171.8335
$ws.Range("N71").Value = -17488

$ws.Range("H110").Value = 47822
$ws.Range("J110").Value = 47822
$ws.Range("L110").Value = 47822
$ws.Range("N110").Value = -56002

$ws.Range("H113").Value = 7410005
$ws.Range("I113").Value = 10103319
$ws.Range("J113").Value = 3392.25
$ws.Range("K113").Value = 10103319
$ws.Range("L113").Value = 3392.25
$ws.Range("M113").Value = -10101149
$ws.Range("N113").Value = -7732.25

$ws.Range("H132").Value = 9777.787
$ws.Range("I132").Value = 10258.128
$ws.Range("J132").Value = 7436.125
$ws.Range("K132").Value = 30774.384
$ws.Range("L132").Value = 22308.375
$ws.Range("M132").Value = -28244.384
$ws.Range("N132").Value = -27368.375

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8334917
$ws.Range("J81").Value = 1791.3334
$ws.Range("L81").Value = 3582.6668
$ws.Range("N81").Value = -5704.6668

$ws.Range("H84").Value = 8334917
$ws.Range("J84").Value = 1791.3334
$ws.Range("L84").Value = 17913.334
$ws.Range("N84").Value = -28521.334

$ws.Range("H100").Value = 2644.7
$ws.Range("I100").Value = 1857.2858
$ws.Range("J100").Value = 4482
$ws.Range("K100").Value = 3714.5716
$ws.Range("L100").Value = 8964
$ws.Range("M100").Value = -3173.5716
$ws.Range("N100").Value = -10046

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 38478390
$ws.Range("I132").Value = 43492636
$ws.Range("K132").Value = 130477908
$ws.Range("M132").Value = -130475378

$ws.Range("H136").Value = 3702.7556
$ws.Range("J136").Value = 2220.75
$ws.Range("L136").Value = 6662.25
$ws.Range("N136").Value = -11762.25
